$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows before the old row 168 (old rows 168-181
# shift down to become 171-184, unchanged).
$ws.Rows.Item(168).EntireRow.Insert()
$ws.Rows.Item(168).EntireRow.Insert()
$ws.Rows.Item(168).EntireRow.Insert()

# New row 168: "Especial" quality lot, Provincia de Melipilla
$ws.Cells.Item(168,1).Value = 10
$ws.Cells.Item(168,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168,3).Value = "La Araucanía"
$ws.Cells.Item(168,4).Value = 44516
$ws.Cells.Item(168,5).Value = 9
$ws.Cells.Item(168,6).Value = "Fruta"
$ws.Cells.Item(168,7).Value = 100101
$ws.Cells.Item(168,8).Value = "Berries"
$ws.Cells.Item(168,9).Value = 100112025
$ws.Cells.Item(168,10).Value = "Frutilla"
$ws.Cells.Item(168,11).Value = "Sin especificar"
$ws.Cells.Item(168,12).Value = "Especial"
$ws.Cells.Item(168,13).Value = 180
$ws.Cells.Item(168,14).Value = 12000
$ws.Cells.Item(168,15).Value = 12000
$ws.Cells.Item(168,16).Value = 12000
$ws.Cells.Item(168,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(168,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(168,19).Value = 1714
$ws.Cells.Item(168,20).Value = 7

# New row 169: "Primera" quality lot, Provincia de Melipilla
$ws.Cells.Item(169,1).Value = 10
$ws.Cells.Item(169,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169,3).Value = "La Araucanía"
$ws.Cells.Item(169,4).Value = 44516
$ws.Cells.Item(169,5).Value = 9
$ws.Cells.Item(169,6).Value = "Fruta"
$ws.Cells.Item(169,7).Value = 100101
$ws.Cells.Item(169,8).Value = "Berries"
$ws.Cells.Item(169,9).Value = 100112025
$ws.Cells.Item(169,10).Value = "Frutilla"
$ws.Cells.Item(169,11).Value = "Sin especificar"
$ws.Cells.Item(169,12).Value = "Primera"
$ws.Cells.Item(169,13).Value = 380
$ws.Cells.Item(169,14).Value = 9000
$ws.Cells.Item(169,15).Value = 9000
$ws.Cells.Item(169,16).Value = 9000
$ws.Cells.Item(169,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(169,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(169,19).Value = 1286
$ws.Cells.Item(169,20).Value = 7

# New row 170: "Segunda" quality lot, Provincia de Melipilla
$ws.Cells.Item(170,1).Value = 10
$ws.Cells.Item(170,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170,3).Value = "La Araucanía"
$ws.Cells.Item(170,4).Value = 44516
$ws.Cells.Item(170,5).Value = 9
$ws.Cells.Item(170,6).Value = "Fruta"
$ws.Cells.Item(170,7).Value = 100101
$ws.Cells.Item(170,8).Value = "Berries"
$ws.Cells.Item(170,9).Value = 100112025
$ws.Cells.Item(170,10).Value = "Frutilla"
$ws.Cells.Item(170,11).Value = "Sin especificar"
$ws.Cells.Item(170,12).Value = "Segunda"
$ws.Cells.Item(170,13).Value = 155
$ws.Cells.Item(170,14).Value = 7000
$ws.Cells.Item(170,15).Value = 7000
$ws.Cells.Item(170,16).Value = 7000
$ws.Cells.Item(170,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(170,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(170,19).Value = 1000
$ws.Cells.Item(170,20).Value = 7

Write-Output "Inserted 3 rows at 168-170; sheet now has"
Write-Output $ws.UsedRange.Rows.Count
